$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "71.930.56"
$ws.Range("E2").Value = "  -0.68%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.996.62"
$ws.Range("E3").Value = "  -1.07%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "542.88"
$ws.Range("E5").Value = "  +4.14%  "

# Row 6 - Solana
$ws.Range("D6").Value = "150.23"
$ws.Range("E6").Value = "  +1.52%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.703"
$ws.Range("E7").Value = "  +12.38%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.746"
$ws.Range("E9").Value = "  +1.06%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  -3.39%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "52.05"
$ws.Range("E11").Value = "  +9.70%  "

# Row 12 - ShibaInu
$ws.Range("D12").Value = "0.0000324"
$ws.Range("E12").Value = "  -3.20%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "10.68"
$ws.Range("E13").Value = "  -2.32%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.631.15"
$ws.Range("E14").Value = "  -1.08%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.987.99"
$ws.Range("E15").Value = "  -1.59%  "

# Row 16 - Uniswap
$ws.Range("D16").Value = "14.10"
$ws.Range("E16").Value = "  -0.79%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "20.49"
$ws.Range("E17").Value = "  -3.63%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -0.24%  "

# Row 19 - Polygon
$ws.Range("E19").Value = "  -1.97%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "71.781.30"
$ws.Range("E20").Value = "  -0.81%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "430.36"
$ws.Range("E21").Value = "  -2.30%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "97.25"
$ws.Range("E22").Value = "  -2.60%  "

# Row 23 - ImmutableX
$ws.Range("D23").Value = "3.52"
$ws.Range("E23").Value = "  -0.81%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  +5.79%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").Value = "14.36"
$ws.Range("E25").Value = "  -2.38%  "

# Row 26 - RenderToken
$ws.Range("D26").Value = "11.51"
$ws.Range("E26").Value = "  -3.42%  "

# Row 27 - Filecoin
$ws.Range("D27").Value = "10.73"
$ws.Range("E27").Value = "  -5.03%  "

# Row 28 - LEO
$ws.Range("D28").Value = "5.86"
$ws.Range("E28").Value = "  +1.04%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "36.77"
$ws.Range("E29").Value = "  -2.55%  "

# Row 30 - Toncoin
$ws.Range("D30").Value = "3.60"

# Row 31 - was Cosmos, now NEARProtocol
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "7.29"
$ws.Range("E31").Value = "  +4.89%  "

# Row 32 - was Hedera, now InjectiveProtocol
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "49.53"
$ws.Range("E32").Value = "  +20.30%  "

# Row 33 - was InjectiveProtocol, now Cosmos
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "13.43"
$ws.Range("E33").Value = "  -0.82%  "

# Row 34 - was NEARProtocol, now Hedera
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.131"
$ws.Range("E34").Value = "  +1.68%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "676.82"
$ws.Range("E35").Value = "  -2.02%  "

# Row 36 - OKB
$ws.Range("D36").Value = "66.03"
$ws.Range("E36").Value = "  -2.93%  "

# Row 37 - TheGraph
$ws.Range("E37").Value = "  +0.49%  "

# Row 38 - PEPE
$ws.Range("E38").Value = "  -7.30%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -1.88%  "

# Row 40 - ThetaToken
$ws.Range("D40").Value = "3.42"
$ws.Range("E40").Value = "  -7.65%  "

# Row 41 - was Dai, now WEMIXToken
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "3.34"
$ws.Range("E41").Value = "  +5.76%  "

# Row 42 - was WEMIXToken, now Dai
$ws.Range("B42").Value = "Dai"
$ws.Range("C42").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.01%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.25%  "

# Row 44 - VeChain
$ws.Range("D44").Value = "0.0488"
$ws.Range("E44").Value = "  -0.85%  "

# Row 45 - Stellar
$ws.Range("D45").Value = "0.150"
$ws.Range("E45").Value = "  +2.21%  "

# Row 46 - Fetch.AI
$ws.Range("D46").Value = "2.73"
$ws.Range("E46").Value = "  -2.28%  "

# Row 47 - THORChain
$ws.Range("E47").Value = "  +8.29%  "

# Row 48 - ApeXProtocol
$ws.Range("D48").Value = "3.37"
$ws.Range("E48").Value = "  -3.78%  "

# Row 49 - FLOKI
$ws.Range("D49").Value = "0.000279"
$ws.Range("E49").Value = "  -0.09%  "

# Row 50 - Stacks
$ws.Range("E50").Value = "  -3.83%  "

# Row 51 - Monero
$ws.Range("D51").Value = "145.01"
$ws.Range("E51").Value = "  +1.65%  "
